$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ext_mgmt_psi* labels in column A (rows 3-5) so the shared
# string text reflects the combined serotype naming scheme.
$ws.Range("A3").Value = "ext_mgmt_psi1_2_3::s1_mean_burden"
$ws.Range("A4").Value = "ext_mgmt_psi1_2_3::s2_mean_burden"
$ws.Range("A5").Value = "ext_mgmt_psi1_2_3::s3_mean_burden"

# Move the window/scroll position so the frozen pane now starts at column B
# instead of H, and select A6 as the active cell in the bottom-right pane.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 3
$win.Left = 12120
$win.Top = 2100

$ws.Range("A6").Select() | Out-Null
